$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the coin-ranking snapshot: update Price (column D) and
# Volume(1h) (column E) for each listed row to the latest scrape.
# Column D holds values like "67.115.86" / "0.0000269" that must stay
# literal text (not be reparsed as numbers, which would also eat
# significant trailing zeros), so force a text format before writing.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.115.86"
$ws.Range("E2").Value = "  -2.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.589.92"
$ws.Range("E3").Value = "  -3.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.48"
$ws.Range("E5").Value = "  -7.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "191.87"
$ws.Range("E6").Value = "  -1.94%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.584.31"
$ws.Range("E7").Value = "  -3.29%  "
$ws.Range("E8").Value = "  -2.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.997"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.680"
$ws.Range("E10").Value = "  -6.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.151"
$ws.Range("E11").Value = "  -6.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "56.14"
$ws.Range("E12").Value = "  -7.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000269"
$ws.Range("E13").Value = "  -6.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.91"
$ws.Range("E14").Value = "  -5.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.155.59"
$ws.Range("E15").Value = "  -3.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.573.37"
$ws.Range("E16").Value = "  -3.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.126"
$ws.Range("E17").Value = "  -1.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.41"
$ws.Range("E18").Value = "  -5.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "67.085.83"
$ws.Range("E19").Value = "  -1.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.21"
$ws.Range("E20").Value = "  -5.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.07"
$ws.Range("E21").Value = "  -7.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "398.51"
$ws.Range("E22").Value = "  -2.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.22"
$ws.Range("E23").Value = "  -10.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.96"
$ws.Range("E24").Value = "  -4.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.24"
$ws.Range("E25").Value = "  -2.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.95"
$ws.Range("E26").Value = "  -4.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.54"
$ws.Range("E27").Value = "  -4.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.07"
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.64"
$ws.Range("E29").Value = "  -3.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.98"
$ws.Range("E30").Value = "  -7.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.72"
$ws.Range("E31").Value = "  +1.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.14"
$ws.Range("E32").Value = "  -5.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "639.32"
$ws.Range("E33").Value = "  +0.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.22"
$ws.Range("E34").Value = "  -3.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.115"
$ws.Range("E35").Value = "  -6.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "63.72"
$ws.Range("E36").Value = "  -6.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.35"
$ws.Range("E37").Value = "  -9.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.406"
$ws.Range("E38").Value = "  -1.77%  "
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0762"
$ws.Range("E40").Value = "  -7.67%  "
$ws.Range("E41").Value = "  -4.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.119.57"
$ws.Range("E42").Value = "  +6.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.72"
$ws.Range("E44").Value = "  +4.05%  "
$ws.Range("E45").Value = "  -3.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0416"
$ws.Range("E46").Value = "  -6.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.131"
$ws.Range("E47").Value = "  -6.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.10"
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.30"
$ws.Range("E49").Value = "  -4.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.54"
$ws.Range("E50").Value = "  -9.92%  "
$ws.Range("E51").Value = "  -0.05%  "
